$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# --- Sheet1 ("data"): append 5 new monthly VIX rows (Feb-Jun 2024) ---
$ws1.Range("A459").Value2 = 45323
$ws1.Range("B459").Value2 = 13.980476190476191

$ws1.Range("A460").Value2 = 45352
$ws1.Range("B460").Value2 = 13.7875

$ws1.Range("A461").Value2 = 45383
$ws1.Range("B461").Value2 = 16.137727272727272

$ws1.Range("A462").Value2 = 45413
$ws1.Range("B462").Value2 = 13.058695652173913

$ws1.Range("A463").Value2 = 45444
$ws1.Range("B463").Value2 = 12.667

# Match formatting of the new rows to the row immediately above them
$ws1.Range("A458:B458").Copy() | Out-Null
$ws1.Range("A459:B463").PasteSpecial(-4122) | Out-Null

# --- Sheet2 ("readme"): refresh the summary labels to reflect the new data range ---
$ws2.Range("B2").Value2 = "to June 2024"
$ws2.Range("A4").Value2 = "Mar 90 - Jun 24"

# --- Restore sheet view / selection state ---
# Set the (inactive) readme sheet's selection first, then activate/select on
# the data sheet last so it ends up as the active tab, matching the source file.
$ws2.Range("A5").Select() | Out-Null

$ws1.Activate()
$ws1.Range("D456").Select() | Out-Null
